# Regenerate save_data K column (column G) with newly computed strikeout
# values (previously derived from "Strike#", now using the pitcher's
# actual strikeout total K). Write the recalculated values in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value for column G ("K")
$newK = @{
    2  = 7
    3  = 4
    4  = 3
    5  = 7
    6  = 6
    7  = 3
    8  = 9
    9  = 3
    10 = 7
    11 = 5
    12 = 7
    13 = 9
    14 = 12
    15 = 4
    16 = 3
    17 = 7
    18 = 5
    19 = 4
    20 = 10
    21 = 3
    22 = 8
    23 = 4
    24 = 6
    25 = 6
    26 = 9
    27 = 2
    28 = 8
    29 = 5
    30 = 5
    31 = 5
    32 = 5
    33 = 9
    34 = 4
    35 = 1
    36 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
